$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (was 45406 / 2024-04-24, now 45436 / 2024-05-24)
$ws.Range("A1").Value = 45436

# Update the price column for rows 27-29 (was 94.3, now 203.5)
$ws.Range("D27").Value = 203.5
$ws.Range("D28").Value = 203.5
$ws.Range("D29").Value = 203.5

# Re-apply all merged cell ranges (refresh the merge collection order)
$ws.Range("A1:D1").UnMerge()
$ws.Range("A1:D1").Merge()
$ws.Range("A9:D9").UnMerge()
$ws.Range("A9:D9").Merge()
$ws.Range("A31:D31").UnMerge()
$ws.Range("A31:D31").Merge()
$ws.Range("B29:C29").UnMerge()
$ws.Range("B29:C29").Merge()
$ws.Range("B26:C26").UnMerge()
$ws.Range("B26:C26").Merge()
$ws.Range("B28:C28").UnMerge()
$ws.Range("B28:C28").Merge()
$ws.Range("A11:D11").UnMerge()
$ws.Range("A11:D11").Merge()
$ws.Range("B27:C27").UnMerge()
$ws.Range("B27:C27").Merge()
$ws.Range("A10:D10").UnMerge()
$ws.Range("A10:D10").Merge()
